# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Each entry gives the target row on the "展览" sheet, the corresponding
# row on the "全部类型" sheet (which contains one extra row not present
# on "展览"), and the new value for column F.
$updates = @(
    @{Sheet1Row = 2;  Sheet4Row = 2;  New = 89},
    @{Sheet1Row = 4;  Sheet4Row = 5;  New = 276},
    @{Sheet1Row = 6;  Sheet4Row = 7;  New = 10161},
    @{Sheet1Row = 8;  Sheet4Row = 9;  New = 924},
    @{Sheet1Row = 9;  Sheet4Row = 10; New = 1263},
    @{Sheet1Row = 10; Sheet4Row = 11; New = 6601},
    @{Sheet1Row = 11; Sheet4Row = 12; New = 14},
    @{Sheet1Row = 12; Sheet4Row = 13; New = 424},
    @{Sheet1Row = 15; Sheet4Row = 16; New = 3130},
    @{Sheet1Row = 17; Sheet4Row = 18; New = 302},
    @{Sheet1Row = 18; Sheet4Row = 19; New = 614},
    @{Sheet1Row = 19; Sheet4Row = 20; New = 118},
    @{Sheet1Row = 21; Sheet4Row = 22; New = 272},
    @{Sheet1Row = 22; Sheet4Row = 23; New = 39},
    @{Sheet1Row = 23; Sheet4Row = 24; New = 1567}
)

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($u in $updates) {
    $ws1.Cells.Item($u.Sheet1Row, 6).Value = $u.New
    $ws4.Cells.Item($u.Sheet4Row, 6).Value = $u.New
}
